{"js": "// Update the lattice-multiplication exercise table: replace the 5-line\n// contents (problem, two multiplier digits, dashed rule, two partial-product\n// placeholder rows) of every cell with a new set of problems/digits.\n//\n// Each table cell holds a single run containing 5 <w:t> text nodes separated\n// by manual line breaks (<w:br/>). We rebuild each cell's paragraph via\n// insertOoxml(..., Replace) so the exact OOXML shape (run properties,\n// <w:br/> placement, and xml:space=\"preserve\" on lines with leading/\n// trailing spaces) is reproduced faithfully.\n\n// New content for every cell, in row-major order (1-based row/col).\nconst newCells = [\n  { row: 1, col: 1, lines: [\"14 x 95\", \"  9    5\", \"  ----\", \"1|    |\", \"4|    |\"] },\n  { row: 1, col: 2, lines: [\"74 x 33\", \"  3    3\", \"  ----\", \"7|    |\", \"4|    |\"] },\n  { row: 1, col: 3, lines: [\"19 x 57\", \"  5    7\", \"  ----\", \"1|    |\", \"9|    |\"] },\n  { row: 2, col: 1, lines: [\"37 x 49\", \"  4    9\", \"  ----\", \"3|    |\", \"7|    |\"] },\n  { row: 2, col: 2, lines: [\"60 x 36\", \"  3    6\", \"  ----\", \"6|    |\", \"0|    |\"] },\n  { row: 2, col: 3, lines: [\"97 x 36\", \"  3    6\", \"  ----\", \"9|    |\", \"7|    |\"] },\n  { row: 3, col: 1, lines: [\"29 x 48\", \"  4    8\", \"  ----\", \"2|    |\", \"9|    |\"] },\n  { row: 3, col: 2, lines: [\"21 x 46\", \"  4    6\", \"  ----\", \"2|    |\", \"1|    |\"] },\n  { row: 3, col: 3, lines: [\"37 x 76\", \"  7    6\", \"  ----\", \"3|    |\", \"7|    |\"] },\n  { row: 4, col: 1, lines: [\"97 x 80\", \"  8    0\", \"  ----\", \"9|    |\", \"7|    |\"] },\n  { row: 4, col: 2, lines: [\"74 x 89\", \"  8    9\", \"  ----\", \"7|    |\", \"4|    |\"] },\n  { row: 4, col: 3, lines: [\"41 x 37\", \"  3    7\", \"  ----\", \"4|    |\", \"1|    |\"] },\n  { row: 5, col: 1, lines: [\"64 x 78\", \"  7    8\", \"  ----\", \"6|    |\", \"4|    |\"] },\n  { row: 5, col: 2, lines: [\"74 x 46\", \"  4    6\", \"  ----\", \"7|    |\", \"4|    |\"] },\n  { row: 5, col: 3, lines: [\"39 x 52\", \"  5    2\", \"  ----\", \"3|    |\", \"9|    |\"] },\n];\n\nfunction xmlEscape(s) {\n  return s\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\")\n    .replace(/\"/g, \"&quot;\");\n}\n\n// Build the <w:t>...</w:t> run contents for one cell's five lines, joined by\n// <w:br/>, adding xml:space=\"preserve\" whenever a line has leading/trailing\n// whitespace (matching how Word itself marks such runs).\nfunction buildRunXml(lines) {\n  return lines\n    .map((line) => {\n      const needsPreserve = /^\\s|\\s$/.test(line);\n      const attr = needsPreserve ? ' xml:space=\"preserve\"' : \"\";\n      return `<w:t${attr}>${xmlEscape(line)}</w:t>`;\n    })\n    .join(\"<w:br/>\");\n}\n\nfunction buildFlatOpc(lines) {\n  const runXml = buildRunXml(lines);\n  return (\n    '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    \"<w:body><w:p><w:r><w:rPr><w:sz w:val=\\\"32\\\"/></w:rPr>\" +\n    runXml +\n    \"</w:r></w:p></w:body></w:document>\" +\n    \"</pkg:xmlData></pkg:part></pkg:package>\"\n  );\n}\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nfor (const spec of newCells) {\n  const cell = table.getCell(spec.row - 1, spec.col - 1);\n  const flatOpc = buildFlatOpc(spec.lines);\n  cell.body.insertOoxml(flatOpc, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the lattice-multiplication exercise table: replace the 5-line\n# contents (problem, two multiplier digits, dashed rule, two partial-product\n# placeholder rows) of every cell with the new set of problems/digits.\n#\n# Each cell's text is 5 \"lines\" joined by manual line breaks (w:br, which\n# the Word object model represents as Chr(11)/vertical-tab inside Range.Text).\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$lb = [char]11\n\n# New content for every cell, in row-major order (row, col, lines[5]).\n$newCells = @(\n    @{ Row = 1; Col = 1; Lines = @(\"14 x 95\", \"  9    5\", \"  ----\", \"1|    |\", \"4|    |\") },\n    @{ Row = 1; Col = 2; Lines = @(\"74 x 33\", \"  3    3\", \"  ----\", \"7|    |\", \"4|    |\") },\n    @{ Row = 1; Col = 3; Lines = @(\"19 x 57\", \"  5    7\", \"  ----\", \"1|    |\", \"9|    |\") },\n    @{ Row = 2; Col = 1; Lines = @(\"37 x 49\", \"  4    9\", \"  ----\", \"3|    |\", \"7|    |\") },\n    @{ Row = 2; Col = 2; Lines = @(\"60 x 36\", \"  3    6\", \"  ----\", \"6|    |\", \"0|    |\") },\n    @{ Row = 2; Col = 3; Lines = @(\"97 x 36\", \"  3    6\", \"  ----\", \"9|    |\", \"7|    |\") },\n    @{ Row = 3; Col = 1; Lines = @(\"29 x 48\", \"  4    8\", \"  ----\", \"2|    |\", \"9|    |\") },\n    @{ Row = 3; Col = 2; Lines = @(\"21 x 46\", \"  4    6\", \"  ----\", \"2|    |\", \"1|    |\") },\n    @{ Row = 3; Col = 3; Lines = @(\"37 x 76\", \"  7    6\", \"  ----\", \"3|    |\", \"7|    |\") },\n    @{ Row = 4; Col = 1; Lines = @(\"97 x 80\", \"  8    0\", \"  ----\", \"9|    |\", \"7|    |\") },\n    @{ Row = 4; Col = 2; Lines = @(\"74 x 89\", \"  8    9\", \"  ----\", \"7|    |\", \"4|    |\") },\n    @{ Row = 4; Col = 3; Lines = @(\"41 x 37\", \"  3    7\", \"  ----\", \"4|    |\", \"1|    |\") },\n    @{ Row = 5; Col = 1; Lines = @(\"64 x 78\", \"  7    8\", \"  ----\", \"6|    |\", \"4|    |\") },\n    @{ Row = 5; Col = 2; Lines = @(\"74 x 46\", \"  4    6\", \"  ----\", \"7|    |\", \"4|    |\") },\n    @{ Row = 5; Col = 3; Lines = @(\"39 x 52\", \"  5    2\", \"  ----\", \"3|    |\", \"9|    |\") }\n)\n\nforeach ($cellSpec in $newCells) {\n    $cell = $t.Cell($cellSpec.Row, $cellSpec.Col)\n    $r = $cell.Range\n    # Trim off the trailing cell-mark / paragraph-mark (last 2 chars) so we\n    # only replace the visible text, preserving the cell's own structure.\n    $r.End = $r.End - 2\n    $r.Text = [string]::Join($lb, $cellSpec.Lines)\n}\n"}
